$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("E44").Select()
